$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 43. This shifts the existing rows 43-87
# down to 44-88 (matching the diff) and extends the used range to T88.
$ws.Rows.Item(43).Insert()

# Populate the newly inserted row 43 with its data.
$ws.Range("A43").Value = 5
$ws.Range("B43").Value = "Macroferia Regional de Talca"
$ws.Range("C43").Value = "Maule"
$ws.Range("D43").Value = 44897
$ws.Range("E43").Value = 7
$ws.Range("F43").Value = "Fruta"
$ws.Range("G43").Value = 100101
$ws.Range("H43").Value = "Berries"
$ws.Range("I43").Value = 100101001
$ws.Range("J43").Value = "Arándano (blue)"
$ws.Range("K43").Value = "Sin especificar"
$ws.Range("L43").Value = "Primera"
$ws.Range("M43").Value = 210
$ws.Range("N43").Value = 3000
$ws.Range("O43").Value = 3000
$ws.Range("P43").Value = 3000
$ws.Range("Q43").Value = "$/bandeja 2 kilos"
$ws.Range("R43").Value = "Provincia de Curicó"
$ws.Range("S43").Value = 1500
$ws.Range("T43").Value = 2
